$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") all share the
# same row layout:
#   - Add a header "Fonte/Tecnologia" in A1 (with the same header style as B1)
#   - Fix diacritics in column A labels (rows 2-12)
#   - Remove the bold/border/alignment style from column A labels (rows 2-12)
# ---------------------------------------------------------------------------

$fixedLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

for ($sheetIdx = 1; $sheetIdx -le 4; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    # New header cell A1, copying the header style already present on B1.
    $ws.Range("B1").Copy($ws.Range("A1"))
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    foreach ($r in $fixedLabels.Keys) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.Value = $fixedLabels[$r]
        $cell.ClearFormats()
        $cell.Value = $fixedLabels[$r]
    }
}

# ---------------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais (MtCO2eq)")
#   - Add header "Período" in A1
#   - Fix diacritics in A2/A3, drop their style
#   - Delete row 4 ("Teto") entirely
# ---------------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy($ws5.Range("A1"))
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").ClearFormats()
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6 ("Custo Total (bilhões de R$)")
#   - Add header "Tipo Expansão" in A1
#   - Rename B1 "Custo" -> "2015"
#   - Fix diacritics in A2/A3, drop their style
#   - Update B2/B3 values
# ---------------------------------------------------------------------------

$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy($ws6.Range("A1"))
$ws6.Range("A1").Value = "Tipo Expansão"

# Reuse the "2015" text cell from another sheet so B1 keeps its text type
# (a plain .Value = "2015" assignment gets auto-coerced to a number).
$wb.Worksheets.Item(1).Range("B1").Copy($ws6.Range("B1"))

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 585

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").ClearFormats()
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
